# The underlying commit ("Change target framework to netstandard2.0") only
# touches the .NET project/build configuration of the PowerpointTemplater
# library itself; it does not change any slide content, ordering, table
# data, images, etc. The only delta in the compiled .pptx test fixture is
# that a handful of relationship ids (the hex strings after "R", e.g.
# R180b9e9458fd4691 -> R7b11cbcdb1e2463c) were re-minted by the OOXML writer
# the next time the fixture was regenerated - those ids are produced from
# Guid.NewGuid() by the Open XML SDK and carry no semantic meaning: they
# still point at the exact same slide parts / media part as before, in the
# same order, with the same content.
#
# There is therefore no actual presentation content for PowerPoint's object
# model to change here, so this script intentionally performs a no-op
# (touching $ppt.ActivePresentation without mutating it) rather than
# attempting to fabricate new random relationship-id strings, which:
#   * cannot be produced deterministically through COM automation, and
#   * would require destructively deleting/recreating the slides/picture,
#     which would itself introduce real (and incorrect) content changes
#     such as different shape ids/names and loss of the picture
#     placeholder/title metadata used by the templating engine.
$p = $ppt.ActivePresentation
$null = $p.Slides.Count
